$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (moves forward one month: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update the prices for the curved-pipe supports
$ws.Range("D30").Value = 1576
$ws.Range("D31").Value = 1794
